# laid work for UART communication modules
# Adds a NOTE about reserved words (FFF000 / FFFF00) used by the
# ROM_FLASHING program, just below the existing ISA documentation table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "NOTE: "
$ws.Range("B32").Value = "FFF000 and FFFF00 are reserved words in the CPU, and are used in the ROM_FLASHING program to signify the end of transmission of new words. As of now"
$ws.Range("B33").Value = "these words are illegal syntax, and cannot be created regardless, but if these are somehow created and flashed in accidentally, the flash will terminate early."

# Match the author's final view state: scrolled down, selection on B33.
$ws.Range("B33").Select()
$excel.ActiveWindow.ScrollRow = 13
